# kredi faiz indirimi yontemi degistirildi
# Update the discounted annual interest rate ("indirimli_yillik_faiz")
# on the "binek" sheet from 8% to 10%, and move the active selection
# to D8 (matching the saved cursor position recorded in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("binek")

# indirimli_yillik_faiz (row 5, col B): 0.08 -> 0.1
$ws.Range("B5").Value = 0.1

# Move / record the active cell selection on the sheet
$ws.Activate() | Out-Null
$ws.Range("D8").Select() | Out-Null
